# Lab01 Review Report — add coding phase defects
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Workbook-level calculation setting (iterative calculation delta)
# ---------------------------------------------------------------------------
$excel.IterativeCalculation = $true
$excel.MaxChange = 0.0001

# ---------------------------------------------------------------------------
# "Coding Phase Defects" sheet — fill in the review header + defect table
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Coding Phase Defects")
$ws.Activate()

# Reviewer block (top-right mini table: Student 3 name + group)
$ws.Range("I3").Value = "Silivăstru Oana"
$ws.Range("J3").Value = 236

# Reviewer name / review date (left block)
$ws.Range("D6").Value = "Silivăstru Oana Maria"
$ws.Range("D7").Value2 = 45736
$ws.Range("D7").NumberFormat = "mm-dd-yy"

# Clear the "Crt. No." header text in B9 (kept blank in the final sheet)
$ws.Range("B9").Value = ""

# Defects table (Crt.No already has B10:B18 autonumbered formulas)
$ws.Range("C10").Value = "C06"
$ws.Range("D10").Value = "LinkedTaskList.java/47"
$ws.Range("E10").Value = "task is an input parameter which is not checked if its valid or not (null)."
$ws.Range("A10").RowHeight = 30

$ws.Range("C11").Value = "C08"
$ws.Range("D11").Value = "Task.java/42"
$ws.Range("E11").Value = 'The condition is (interval < 1) and the message says "interval should me > 1". There is a typo and an error, it should say "interval should be > 0"'
$ws.Range("A11").RowHeight = 60

$ws.Range("C12").Value = "C01"
$ws.Range("D12").Value = "Task.java/107"
$ws.Range("E12").Value = "This condition is all the time true because all of its check components have already been checked before. It introduces confusion."
$ws.Range("A12").RowHeight = 15

$ws.Range("C13").Value = "C11"
$ws.Range("D13").Value = "LinkedTaskList.java/98"
$ws.Range("E13").Value = "The variable name is tks instead of tasks. At first view, anyone who read the code might get confused because tks does not really have a meaning."
$ws.Range("A13").RowHeight = 60

$ws.Range("C14").Value = "C01"
$ws.Range("D14").Value = "LinkedTaskList.java/15"
$ws.Range("E14").Value = "cursor variable is not initialized. Initialized it to 0."
$ws.Range("A14").RowHeight = 30

$ws.Range("C15").Value = "C07"
$ws.Range("D15").Value = "LinkedTaskList.java/39"
$ws.Range("E15").Value = "cursor variable now equals lastCalled which is the index of the removed element. What will happen if we call remove twice? Instead, if it's not the first element, then after removal just decrease the cursor index."
$ws.Range("A15").RowHeight = 75

$ws.Range("C16").Value = "C06"
$ws.Range("D16").Value = "LinkedTaskList.java/60"
$ws.Range("E16").Value = "this.last is not verified to not be null (empty list) . Added a check for that."
$ws.Range("A16").RowHeight = 30

$ws.Range("C17").Value = "C08"
$ws.Range("D17").Value = "Task.java/25"
$ws.Range("E17").Value = "Constructor does not have the thorws keyword in the definition even if it can throw an exception. Added it."
$ws.Range("A17").RowHeight = 45

$ws.Range("C18").Value = "C01"
$ws.Range("D18").Value = "Task.java/134"
$ws.Range("E18").Value = "No need for else; either the function enters the if and then returns or just returns. Removed else."
$ws.Range("A18").RowHeight = 45

# Effort to review document
$ws.Range("E32").Value = "0.5h"

# View state: scroll/zoom/selection for this sheet
$ws.Range("J3").Select()
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1

# ---------------------------------------------------------------------------
# Other sheets — view state only (zoom / scroll / selection)
# ---------------------------------------------------------------------------
$wsReq = $wb.Worksheets.Item("Requirements Phase Defects")
$wsReq.Activate()
$wsReq.Range("E10").Select()
$excel.ActiveWindow.Zoom = 145
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1

$wsArch = $wb.Worksheets.Item("Architect. Design Phase Defects")
$wsArch.Activate()
$wsArch.Range("I3").Select()

$wsTool = $wb.Worksheets.Item("Tool-basedCodeAnalysis")
$wsTool.Activate()
$wsTool.Range("F36").Select()

# Re-activate the Coding Phase Defects sheet last (it is the active tab)
$ws.Activate()
$wb.Windows.Item(1).DisplayWorkbookTabs = $true
